# Updates to v 0.32
# - Corrects the NV distance label from 0.45 m to 0.40 m on the
#   visual-assistance slides (5-8).
# - Refreshes the cached "datetimeFigureOut" footer field (3/21/24 -> 4/9/24)
#   on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# --- 1) Fix the NV distance text boxes on slides 5-8 -----------------------
$nvDistanceShapeNames = @{
    5 = "TextBox 13"
    6 = "TextBox 11"
    7 = "TextBox 10"
    8 = "TextBox 10"
}

foreach ($slideIndex in $nvDistanceShapeNames.Keys) {
    $slide = $p.Slides.Item($slideIndex)
    $shapeName = $nvDistanceShapeNames[$slideIndex]
    $shape = $slide.Shapes.Item($shapeName)
    if ($shape.TextFrame.TextRange.Text -eq "0.45 m") {
        $shape.TextFrame.TextRange.Text = "0.40 m"
    }
}

# --- 2) Refresh the cached date field text everywhere it appears -----------
$oldDate = "3/21/24"
$newDate = "4/9/24"

# Slide master
$masterDateShape = $p.SlideMaster.Shapes.Item("Date Placeholder 3")
if ($masterDateShape.TextFrame.TextRange.Text -eq $oldDate) {
    $masterDateShape.TextFrame.TextRange.Text = $newDate
}

# Every slide layout under the master
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shape = $layout.Shapes.Item($j)
        if ($shape.Name -like "Date Placeholder*") {
            if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
